$d = $word.ActiveDocument

# The document currently has a collapsed "_GoBack" bookmark sitting right
# after the "Metode Regula Falsi" run (near the end of the 2nd bullet of
# the first list). Word moves this special bookmark to track the location
# of the user's last edit, so here it needs to move to the very start of
# the document (immediately after the pPr of the first paragraph, i.e.
# right before the very first run "TUGAS FINAL - ...").

# 1) Remove the old "_GoBack" bookmark from its current location.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 2) Insert a one-character placeholder at the very beginning of the
#    document. (A zero-length Range(0,0) passed straight to
#    Bookmarks.Add collapses incorrectly when it sits at the absolute
#    start of the story, so we briefly widen the insertion point with a
#    throwaway character, bookmark right after it, then delete the
#    character again — leaving a clean collapsed bookmark at position 0.)
$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertBefore("X")

# 3) Add the "_GoBack" bookmark immediately after the placeholder
#    character (a true collapsed range).
$afterPlaceholder = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $afterPlaceholder)

# 4) Delete the placeholder character; the bookmark collapses down to the
#    very start of the document, right after the first paragraph's pPr
#    and before its first run.
$placeholder = $d.Range(0, 1)
$placeholder.Delete()
